$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.051.87'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '3.533.43'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '604.02'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').Value = '143.87'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('D7').Value = '3.531.20'
$ws.Range('E7').Value = '  -1.79%  '
$ws.Range('D9').Value = '0.513'
$ws.Range('E9').Value = '  +4.82%  '
$ws.Range('D10').Value = '7.83'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('E11').Value = '  -4.60%  '
$ws.Range('D12').Value = '0.407'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '4.134.93'
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').Value = '0.0000194'
$ws.Range('E14').Value = '  -7.48%  '
$ws.Range('D15').Value = '28.34'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('D16').Value = '3.525.88'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '65.922.24'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').Value = '10.93'
$ws.Range('E19').Value = '  -5.17%  '
$ws.Range('D20').Value = '6.19'
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('D21').Value = '14.53'
$ws.Range('E21').Value = '  -3.61%  '
$ws.Range('D22').Value = '421.84'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').Value = '0.594'
$ws.Range('E23').Value = '  -4.49%  '
$ws.Range('D24').Value = '77.01'
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('D25').Value = '3.675.42'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').Value = '  -6.32%  '
$ws.Range('E28').Value = '  -2.26%  '
$ws.Range('D29').Value = '7.80'
$ws.Range('E29').Value = '  -6.14%  '
$ws.Range('D30').Value = '8.88'
$ws.Range('E30').Value = '  -4.64%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '3.541.51'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').Value = '0.155'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').Value = '24.22'
$ws.Range('E34').Value = '  -5.17%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -7.61%  '
$ws.Range('D37').Value = '7.58'
$ws.Range('E37').Value = '  -3.53%  '
$ws.Range('D38').Value = '177.88'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').Value = '1.63'
$ws.Range('E39').Value = '  -5.53%  '
$ws.Range('D40').Value = '5.21'
$ws.Range('E40').Value = '  -7.41%  '
$ws.Range('D41').Value = '0.0816'
$ws.Range('E41').Value = '  -5.07%  '
$ws.Range('D42').Value = '0.860'
$ws.Range('E42').Value = '  -4.23%  '
$ws.Range('D43').Value = '4.97'
$ws.Range('E43').Value = '  -5.40%  '
$ws.Range('D44').Value = '45.49'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('E45').Value = '  -7.71%  '
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '23.66'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  -8.64%  '
$ws.Range('D49').Value = '7.04'
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('D50').Value = '1.12'
$ws.Range('E50').Value = '  -6.51%  '
$ws.Range('D51').Value = '0.904'
$ws.Range('E51').Value = '  -5.20%  '
